# Auto-generated Excel COM-interop script replicating the cryptos.xlsx refresh commit.
# Updates the "Price" (D) and "Volume(1h)" (E) columns for each coin row, and
# re-labels rows 42-46 whose ranking order changed (coin identity + link move rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while guaranteeing it stays plain TEXT
# (matches the source data, which stores every Price/Volume cell as a string,
# even when the text looks like a plain number, e.g. "0.999" or "0.0000154").
# Temporarily forcing the Text number format keeps Excel from auto-converting
# the input to a numeric type; resetting the Style back to "Normal" afterwards
# avoids leaving a stray number-format override behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Rows 2-41 & 47-51: Price / Volume(1h) refresh ---
$ws.Range("D2").Value = "64.117.96"
$ws.Range("E2").Value = "  +5.30%  "
$ws.Range("D3").Value = "2.731.18"
$ws.Range("E3").Value = "  +3.25%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.15%  "
Set-TextValue $ws.Range("D5") "581.42"
$ws.Range("E5").Value = "  +0.26%  "
Set-TextValue $ws.Range("D6") "157.30"
$ws.Range("E6").Value = "  +9.18%  "
Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.14%  "
Set-TextValue $ws.Range("D8") "0.617"
$ws.Range("E8").Value = "  +3.33%  "
$ws.Range("D9").Value = "2.755.93"
$ws.Range("E9").Value = "  +3.53%  "
Set-TextValue $ws.Range("D10") "6.78"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("E11").Value = "  +3.74%  "
Set-TextValue $ws.Range("D12") "0.393"
$ws.Range("E12").Value = "  +4.07%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "3.220.96"
$ws.Range("E14").Value = "  +3.21%  "
Set-TextValue $ws.Range("D15") "27.52"
$ws.Range("E15").Value = "  +4.24%  "
$ws.Range("D16").Value = "63.740.60"
$ws.Range("E16").Value = "  +4.76%  "
Set-TextValue $ws.Range("D17") "0.0000154"
$ws.Range("E17").Value = "  +7.17%  "
$ws.Range("D18").Value = "2.749.84"
$ws.Range("E18").Value = "  +3.62%  "
Set-TextValue $ws.Range("D19") "12.06"
$ws.Range("E19").Value = "  +3.58%  "
Set-TextValue $ws.Range("D20") "4.94"
$ws.Range("E20").Value = "  +4.21%  "
Set-TextValue $ws.Range("D21") "362.53"
$ws.Range("E21").Value = "  +2.98%  "
Set-TextValue $ws.Range("D22") "6.95"
$ws.Range("E22").Value = "  +0.72%  "
Set-TextValue $ws.Range("D23") "0.544"
$ws.Range("E23").Value = "  +3.04%  "
Set-TextValue $ws.Range("D24") "0.997"
$ws.Range("E24").Value = "  -0.26%  "
Set-TextValue $ws.Range("D25") "66.78"
$ws.Range("E25").Value = "  +4.48%  "
$ws.Range("E26").Value = "  +5.51%  "
Set-TextValue $ws.Range("D27") "8.64"
$ws.Range("E27").Value = "  +2.92%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "0.0₃0917"
$ws.Range("E29").Value = "  +13.33%  "
Set-TextValue $ws.Range("D30") "2.04"
$ws.Range("E30").Value = "  +1.65%  "
Set-TextValue $ws.Range("D31") "7.22"
$ws.Range("E31").Value = "  +6.31%  "
Set-TextValue $ws.Range("D32") "1.28"
$ws.Range("E32").Value = "  +19.24%  "
Set-TextValue $ws.Range("D33") "174.04"
$ws.Range("E33").Value = "  +4.54%  "
$ws.Range("E34").Value = "  -0.08%  "
Set-TextValue $ws.Range("D35") "20.62"
$ws.Range("E35").Value = "  +3.10%  "
Set-TextValue $ws.Range("D36") "4.91"
$ws.Range("E36").Value = "  +7.42%  "
Set-TextValue $ws.Range("D37") "1.47"
$ws.Range("E37").Value = "  +10.49%  "
Set-TextValue $ws.Range("D38") "1.83"
$ws.Range("E38").Value = "  +8.39%  "
Set-TextValue $ws.Range("D39") "1.01"
$ws.Range("E39").Value = "  +11.76%  "
Set-TextValue $ws.Range("D40") "4.29"
$ws.Range("E40").Value = "  +3.73%  "
Set-TextValue $ws.Range("D41") "338.09"
$ws.Range("E41").Value = "  -1.32%  "
Set-TextValue $ws.Range("D47") "0.648"
$ws.Range("E47").Value = "  +3.90%  "
Set-TextValue $ws.Range("D48") "0.0261"
$ws.Range("E48").Value = "  +4.40%  "
Set-TextValue $ws.Range("D49") "138.93"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("E50").Value = "  +2.92%  "
Set-TextValue $ws.Range("D51") "0.998"
$ws.Range("E51").Value = "  +0.00%  "

# --- Rows 42-46: coin ranking reshuffled (RenderToken/OKB and Hedera/InjectiveProtocol swap places) ---
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D42") "6.02"
$ws.Range("E42").Value = "  +15.60%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D43") "39.50"
$ws.Range("E43").Value = "  +3.03%  "

Set-TextValue $ws.Range("D44") "21.92"
$ws.Range("E44").Value = "  +7.82%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D45") "22.23"
$ws.Range("E45").Value = "  +7.55%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D46") "0.0605"
$ws.Range("E46").Value = "  +5.38%  "

